# Remove the closure date row for 04/09/2024 (serial 45539) from Foglio1.
# Deleting the entire row shifts all subsequent dates up by one row and
# drops the now-unused last row from the used range.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

$ws.Rows.Item(28).Delete()

# Update the active selection as recorded after the edit.
$ws.Range("C21").Select()
